$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Col18a1"
$ws.Cells.Item(2,3).Value = "Gpc1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 10.92859066666667
$ws.Cells.Item(2,8).Value = 32.785772
$ws.Cells.Item(2,9).Value = 0.2185558471001832
$ws.Cells.Item(2,10).Value = 0.2185558471001832
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 2.127396333333333
$ws.Cells.Item(2,14).Value = 6.382189
$ws.Cells.Item(2,15).Value = 0.06137654768277986
$ws.Cells.Item(2,16).Value = 0.06137654768277986
$ws.Cells.Item(2,17).Value = 23.24944371276756
$ws.Cells.Item(2,18).Value = 209.244993414908
$ws.Cells.Item(2,19).Value = 0.01341420337089474
$ws.Cells.Item(2,20).Value = 0.01341420337089474

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Col18a1"
$ws.Cells.Item(3,3).Value = "Gpc1"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 10.92859066666667
$ws.Cells.Item(3,8).Value = 32.785772
$ws.Cells.Item(3,9).Value = 0.2185558471001832
$ws.Cells.Item(3,10).Value = 0.2185558471001832
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 3.721182333333333
$ws.Cells.Item(3,14).Value = 11.163547
$ws.Cells.Item(3,15).Value = 0.1073581454191429
$ws.Cells.Item(3,16).Value = 0.1073581454191429
$ws.Cells.Item(3,17).Value = 40.66727851703155
$ws.Cells.Item(3,18).Value = 366.005506653284
$ws.Cells.Item(3,19).Value = 0.02346375041518542
$ws.Cells.Item(3,20).Value = 0.02346375041518542

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Col18a1"
$ws.Cells.Item(4,3).Value = "Gpc1"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 10.92859066666667
$ws.Cells.Item(4,8).Value = 32.785772
$ws.Cells.Item(4,9).Value = 0.2185558471001832
$ws.Cells.Item(4,10).Value = 0.2185558471001832
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 28.81280933333333
$ws.Cells.Item(4,14).Value = 86.43842799999999
$ws.Cells.Item(4,15).Value = 0.8312653068980773
$ws.Cells.Item(4,16).Value = 0.8312653068980772
$ws.Cells.Item(4,17).Value = 314.8833991607128
$ws.Cells.Item(4,18).Value = 2833.950592446416
$ws.Cells.Item(4,19).Value = 0.181677893314103
$ws.Cells.Item(4,20).Value = 0.181677893314103

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Col18a1"
$ws.Cells.Item(5,3).Value = "Gpc1"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 21.275312
$ws.Cells.Item(5,8).Value = 63.825936
$ws.Cells.Item(5,9).Value = 0.4254751576214852
$ws.Cells.Item(5,10).Value = 0.4254751576214852
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 2.127396333333333
$ws.Cells.Item(5,14).Value = 6.382189
$ws.Cells.Item(5,15).Value = 0.06137654768277986
$ws.Cells.Item(5,16).Value = 0.06137654768277986
$ws.Cells.Item(5,17).Value = 45.26102073932267
$ws.Cells.Item(5,18).Value = 407.349186653904
$ws.Cells.Item(5,19).Value = 0.02611419629959336
$ws.Cells.Item(5,20).Value = 0.02611419629959336

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Col18a1"
$ws.Cells.Item(6,3).Value = "Gpc1"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 21.275312
$ws.Cells.Item(6,8).Value = 63.825936
$ws.Cells.Item(6,9).Value = 0.4254751576214852
$ws.Cells.Item(6,10).Value = 0.4254751576214852
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 3.721182333333333
$ws.Cells.Item(6,14).Value = 11.163547
$ws.Cells.Item(6,15).Value = 0.1073581454191429
$ws.Cells.Item(6,16).Value = 0.1073581454191429
$ws.Cells.Item(6,17).Value = 79.16931515055467
$ws.Cells.Item(6,18).Value = 712.523836354992
$ws.Cells.Item(6,19).Value = 0.04567822384416014
$ws.Cells.Item(6,20).Value = 0.04567822384416014

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Col18a1"
$ws.Cells.Item(7,3).Value = "Gpc1"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 21.275312
$ws.Cells.Item(7,8).Value = 63.825936
$ws.Cells.Item(7,9).Value = 0.4254751576214852
$ws.Cells.Item(7,10).Value = 0.4254751576214852
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 28.81280933333333
$ws.Cells.Item(7,14).Value = 86.43842799999999
$ws.Cells.Item(7,15).Value = 0.8312653068980773
$ws.Cells.Item(7,16).Value = 0.8312653068980772
$ws.Cells.Item(7,17).Value = 613.0015081631786
$ws.Cells.Item(7,18).Value = 5517.013573468607
$ws.Cells.Item(7,19).Value = 0.3536827374777317
$ws.Cells.Item(7,20).Value = 0.3536827374777317

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Col18a1"
$ws.Cells.Item(8,3).Value = "Gpc1"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 17.79975
$ws.Cells.Item(8,8).Value = 53.39925
$ws.Cells.Item(8,9).Value = 0.3559689952783316
$ws.Cells.Item(8,10).Value = 0.3559689952783316
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 2.127396333333333
$ws.Cells.Item(8,14).Value = 6.382189
$ws.Cells.Item(8,15).Value = 0.06137654768277986
$ws.Cells.Item(8,16).Value = 0.06137654768277986
$ws.Cells.Item(8,17).Value = 37.86712288425
$ws.Cells.Item(8,18).Value = 340.80410595825
$ws.Cells.Item(8,19).Value = 0.02184814801229176
$ws.Cells.Item(8,20).Value = 0.02184814801229176

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Col18a1"
$ws.Cells.Item(9,3).Value = "Gpc1"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 17.79975
$ws.Cells.Item(9,8).Value = 53.39925
$ws.Cells.Item(9,9).Value = 0.3559689952783316
$ws.Cells.Item(9,10).Value = 0.3559689952783316
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 3.721182333333333
$ws.Cells.Item(9,14).Value = 11.163547
$ws.Cells.Item(9,15).Value = 0.1073581454191429
$ws.Cells.Item(9,16).Value = 0.1073581454191429
$ws.Cells.Item(9,17).Value = 66.23611523775
$ws.Cells.Item(9,18).Value = 596.12503713975
$ws.Cells.Item(9,19).Value = 0.03821617115979731
$ws.Cells.Item(9,20).Value = 0.03821617115979731

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Col18a1"
$ws.Cells.Item(10,3).Value = "Gpc1"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 17.79975
$ws.Cells.Item(10,8).Value = 53.39925
$ws.Cells.Item(10,9).Value = 0.3559689952783316
$ws.Cells.Item(10,10).Value = 0.3559689952783316
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 28.81280933333333
$ws.Cells.Item(10,14).Value = 86.43842799999999
$ws.Cells.Item(10,15).Value = 0.8312653068980773
$ws.Cells.Item(10,16).Value = 0.8312653068980772
$ws.Cells.Item(10,17).Value = 512.860802931
$ws.Cells.Item(10,18).Value = 4615.747226379
$ws.Cells.Item(10,19).Value = 0.2959046761062425
$ws.Cells.Item(10,20).Value = 0.2959046761062425
